$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("H1").Value = "createdAt"
$ws.Range("I1").Value = "updatedAt"

$timestamps = @(
    "2023-04-16T03:42:46.836Z",
    "2023-04-16T17:39:00.184Z",
    "2023-04-15T06:41:48.371Z",
    "2023-04-20T12:35:38.970Z",
    "2023-04-12T17:49:56.026Z",
    "2023-04-20T15:33:32.873Z",
    "2023-04-27T17:38:49.996Z",
    "2023-04-21T01:17:39.404Z",
    "2023-04-24T02:45:36.900Z",
    "2023-04-21T22:05:53.134Z",
    "2023-04-23T05:50:33.673Z",
    "2023-04-22T06:35:34.274Z",
    "2023-04-28T06:11:34.174Z",
    "2023-05-01T02:48:49.009Z",
    "2023-04-28T07:29:25.785Z",
    "2023-04-26T10:52:05.915Z",
    "2023-05-07T16:06:05.180Z",
    "2023-05-02T23:32:48.918Z",
    "2023-05-02T02:18:11.393Z",
    "2023-05-06T18:22:45.620Z",
    "2023-05-09T08:22:07.674Z",
    "2023-05-04T20:26:24.705Z",
    "2023-05-11T08:49:28.642Z",
    "2023-05-04T07:50:49.807Z",
    "2023-05-09T12:30:46.024Z",
    "2023-05-12T19:30:19.866Z",
    "2023-05-09T14:25:42.860Z",
    "2023-05-04T09:09:48.198Z",
    "2023-05-03T19:35:51.139Z"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $timestamps[$i]
    $ws.Cells.Item($row, 9).Value = $timestamps[$i]
}

# Update sheet view: clear frozen/scrolled top-left cell and change selection
$ws.Range("F2").Select() | Out-Null
